$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new value. Values are written with a
# leading apostrophe so Excel stores them as literal text (matching the source
# data, which never holds true numbers/percentages -- e.g. "36.921.42" and
# "  +1.83%  " -- and the style is reset afterwards so the quote-prefix does not
# leave a lingering number-format change on the cell.
$updates = @(
    @{ Cell = "D2"; Value = '36.921.42' }
    @{ Cell = "E2"; Value = '  +1.83%  ' }
    @{ Cell = "D3"; Value = '2.033.07' }
    @{ Cell = "E3"; Value = '  +1.01%  ' }
    @{ Cell = "E4"; Value = '  +0.01%  ' }
    @{ Cell = "D5"; Value = '249.26' }
    @{ Cell = "E5"; Value = '  -1.09%  ' }
    @{ Cell = "E6"; Value = '  -0.40%  ' }
    @{ Cell = "D7"; Value = '63.21' }
    @{ Cell = "E7"; Value = '  +0.59%  ' }
    @{ Cell = "D9"; Value = '0.397' }
    @{ Cell = "E9"; Value = '  +7.39%  ' }
    @{ Cell = "D10"; Value = '58.17' }
    @{ Cell = "E10"; Value = '  -1.57%  ' }
    @{ Cell = "D11"; Value = '0.0792' }
    @{ Cell = "E11"; Value = '  +6.38%  ' }
    @{ Cell = "E12"; Value = '  -0.51%  ' }
    @{ Cell = "B13"; Value = 'Polygon' }
    @{ Cell = "C13"; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Cell = "D13"; Value = '0.895' }
    @{ Cell = "E13"; Value = '  -1.10%  ' }
    @{ Cell = "B14"; Value = 'Avalanche' }
    @{ Cell = "C14"; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' }
    @{ Cell = "D14"; Value = '23.45' }
    @{ Cell = "E14"; Value = '  +18.33%  ' }
    @{ Cell = "D15"; Value = '14.42' }
    @{ Cell = "E15"; Value = '  -2.43%  ' }
    @{ Cell = "D16"; Value = '2.330.36' }
    @{ Cell = "E16"; Value = '  +1.02%  ' }
    @{ Cell = "E17"; Value = '  +3.35%  ' }
    @{ Cell = "D18"; Value = '2.029.19' }
    @{ Cell = "E18"; Value = '  +0.89%  ' }
    @{ Cell = "D19"; Value = '36.836.82' }
    @{ Cell = "E19"; Value = '  +1.81%  ' }
    @{ Cell = "D20"; Value = '72.58' }
    @{ Cell = "E20"; Value = '  +0.84%  ' }
    @{ Cell = "D21"; Value = '0.0₃0883' }
    @{ Cell = "E21"; Value = '  +2.90%  ' }
    @{ Cell = "E22"; Value = '  +2.68%  ' }
    @{ Cell = "D23"; Value = '237.10' }
    @{ Cell = "E23"; Value = '  +1.33%  ' }
    @{ Cell = "E24"; Value = '  +0.14%  ' }
    @{ Cell = "D25"; Value = '2.52' }
    @{ Cell = "E25"; Value = '  -6.13%  ' }
    @{ Cell = "E26"; Value = '  +1.88%  ' }
    @{ Cell = "E27"; Value = '  +4.16%  ' }
    @{ Cell = "E28"; Value = '  +19.59%  ' }
    @{ Cell = "D29"; Value = '160.73' }
    @{ Cell = "E29"; Value = '  -1.37%  ' }
    @{ Cell = "E30"; Value = '  +3.82%  ' }
    @{ Cell = "D31"; Value = '0.122' }
    @{ Cell = "E31"; Value = '  +1.02%  ' }
    @{ Cell = "E32"; Value = '  +1.22%  ' }
    @{ Cell = "E33"; Value = '  -1.27%  ' }
    @{ Cell = "D34"; Value = '0.0627' }
    @{ Cell = "E34"; Value = '  +3.54%  ' }
    @{ Cell = "D35"; Value = '4.55' }
    @{ Cell = "E35"; Value = '  +0.23%  ' }
    @{ Cell = "D36"; Value = '6.63' }
    @{ Cell = "E36"; Value = '  +12.22%  ' }
    @{ Cell = "D37"; Value = '2.39' }
    @{ Cell = "E38"; Value = '  +0.07%  ' }
    @{ Cell = "D40"; Value = '3.21' }
    @{ Cell = "E40"; Value = '  +25.22%  ' }
    @{ Cell = "E41"; Value = '  -0.66%  ' }
    @{ Cell = "E42"; Value = '  +3.40%  ' }
    @{ Cell = "E43"; Value = '  +0.95%  ' }
    @{ Cell = "D45"; Value = '0.0217' }
    @{ Cell = "E45"; Value = '  +0.33%  ' }
    @{ Cell = "D46"; Value = '17.01' }
    @{ Cell = "E46"; Value = '  +2.42%  ' }
    @{ Cell = "D47"; Value = '94.23' }
    @{ Cell = "E47"; Value = '  +0.12%  ' }
    @{ Cell = "E48"; Value = '  -1.77%  ' }
    @{ Cell = "D49"; Value = '1.365.90' }
    @{ Cell = "E49"; Value = '  -4.11%  ' }
    @{ Cell = "E50"; Value = '  -0.30%  ' }
    @{ Cell = "D51"; Value = '2.221.94' }
    @{ Cell = "E51"; Value = '  +1.17%  ' }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.Value = "'" + $u.Value
    $c.Style = 'Normal'
}
